{"js": "// Replace the date title and the 25 math-fact cells per the diff.\n// Several cells share identical original text (e.g. \"47\u00f75=\" appears\n// twice), so we replace hits strictly in document order rather than\n// relying on a single global find/replace per string.\n\nconst replacements = [\n  [\"2024-11-17 Sunday\", [\"2024-11-18 Monday\"]],\n  [\"37\u00f76=\", [\"53\u00f76=\"]],\n  [\"75\u00f73=\", [\"85\u00f75=\"]],\n  [\"84\u00f72=\", [\"49\u00f75=\"]],\n  [\"36\u00f72=\", [\"35\u00f78=\"]],\n  [\"86\u00f74=\", [\"84\u00f76=\"]],\n  [\"47\u00f75=\", [\"21\u00f72=\", \"62\u00f78=\"]], // two occurrences, in document order\n  [\"44\u00f73=\", [\"69\u00f75=\"]],\n  [\"91\u00f78=\", [\"52\u00f76=\"]],\n  [\"12\u00f74=\", [\"80\u00f77=\"]],\n  [\"22\u00f78=\", [\"32\u00f79=\"]],\n  [\"49\u00f76=\", [\"73\u00f74=\"]],\n  [\"75\u00f74=\", [\"89\u00f77=\"]],\n  [\"48\u00f72=\", [\"54\u00f73=\"]],\n  [\"94\u00f77=\", [\"38\u00f75=\"]],\n  [\"89\u00f72=\", [\"99\u00f79=\"]],\n  [\"83\u00f78=\", [\"85\u00f79=\"]],\n  [\"10\u00f78=\", [\"93\u00f75=\"]],\n  [\"46\u00f75=\", [\"75\u00f72=\"]],\n  [\"72\u00f79=\", [\"48\u00f75=\"]],\n  [\"59\u00f73=\", [\"14\u00f75=\"]],\n  [\"81\u00f79=\", [\"61\u00f73=\"]],\n  [\"97\u00f75=\", [\"80\u00f74=\"]],\n  [\"26\u00f75=\", [\"11\u00f78=\"]],\n  [\"49\u00f73=\", [\"58\u00f77=\"]],\n];\n\nfor (const [oldText, newTexts] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length !== newTexts.length) {\n    throw new Error(\n      `Expected ${newTexts.length} match(es) for \"${oldText}\" but found ${results.items.length}`\n    );\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newTexts[i], \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date title and the 25 math-fact table cells per the diff.\n# \"47\u00f75=\" occurs twice with two different replacement targets, so each\n# pair below is applied with wdReplaceOne (ReplaceAll would make both\n# identical). A fresh Find.Execute always starts scanning from the top\n# of $d.Content, so running the two wdReplaceOne calls in document\n# order turns the first remaining match into the right value each time.\n\n$d = $word.ActiveDocument\n\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\n$replacements = @(\n    @(\"2024-11-17 Sunday\", \"2024-11-18 Monday\"),\n    @(\"37\u00f76=\", \"53\u00f76=\"),\n    @(\"75\u00f73=\", \"85\u00f75=\"),\n    @(\"84\u00f72=\", \"49\u00f75=\"),\n    @(\"36\u00f72=\", \"35\u00f78=\"),\n    @(\"86\u00f74=\", \"84\u00f76=\"),\n    @(\"47\u00f75=\", \"21\u00f72=\"),\n    @(\"44\u00f73=\", \"69\u00f75=\"),\n    @(\"91\u00f78=\", \"52\u00f76=\"),\n    @(\"12\u00f74=\", \"80\u00f77=\"),\n    @(\"22\u00f78=\", \"32\u00f79=\"),\n    @(\"49\u00f76=\", \"73\u00f74=\"),\n    @(\"75\u00f74=\", \"89\u00f77=\"),\n    @(\"48\u00f72=\", \"54\u00f73=\"),\n    @(\"94\u00f77=\", \"38\u00f75=\"),\n    @(\"89\u00f72=\", \"99\u00f79=\"),\n    @(\"83\u00f78=\", \"85\u00f79=\"),\n    @(\"47\u00f75=\", \"62\u00f78=\"),\n    @(\"10\u00f78=\", \"93\u00f75=\"),\n    @(\"46\u00f75=\", \"75\u00f72=\"),\n    @(\"72\u00f79=\", \"48\u00f75=\"),\n    @(\"59\u00f73=\", \"14\u00f75=\"),\n    @(\"81\u00f79=\", \"61\u00f73=\"),\n    @(\"97\u00f75=\", \"80\u00f74=\"),\n    @(\"26\u00f75=\", \"11\u00f78=\"),\n    @(\"49\u00f73=\", \"58\u00f77=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $found = $find.Execute($oldText, $false, $false, $false, $false, $false, $true, $wdFindContinue, $false, $newText, $wdReplaceOne)\n\n    if (-not $found) {\n        throw \"Could not find text to replace: $oldText\"\n    }\n}\n"}
